# Auto-generated edit script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.599.38"
$ws.Range("E2").Value = "  +1.80%  "
$ws.Range("D3").Value = "2.169.94"
$ws.Range("E3").Value = "  +3.65%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'229.64"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.44%  "
$ws.Range("E6").Value = "  +1.27%  "
$ws.Range("D7").Value = "'63.27"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +4.80%  "
$ws.Range("E9").Value = "  +2.94%  "
$ws.Range("D10").Value = "'0.0862"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.70%  "
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("E12").Value = "  +7.78%  "
$ws.Range("D13").Value = "2.490.31"
$ws.Range("E13").Value = "  +3.61%  "
$ws.Range("D14").Value = "'22.33"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.60%  "
$ws.Range("E15").Value = "  +3.21%  "
$ws.Range("E16").Value = "  +2.24%  "
$ws.Range("D17").Value = "2.172.14"
$ws.Range("E17").Value = "  +3.69%  "
$ws.Range("D18").Value = "39.594.56"
$ws.Range("E18").Value = "  +2.03%  "
$ws.Range("D19").Value = "'72.57"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.37%  "
$ws.Range("E20").Value = "  +2.05%  "
$ws.Range("E21").Value = "  +1.98%  "
$ws.Range("D22").Value = "'229.04"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.71%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  -1.35%  "
$ws.Range("E25").Value = "  +1.29%  "
$ws.Range("D26").Value = "'9.79"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.94%  "
$ws.Range("D27").Value = "'172.82"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.12%  "
$ws.Range("E28").Value = "  -0.50%  "
$ws.Range("E29").Value = "  -3.51%  "
$ws.Range("E30").Value = "  +2.73%  "
$ws.Range("E31").Value = "  +8.57%  "
$ws.Range("E32").Value = "  +1.42%  "
$ws.Range("D33").Value = "'4.68"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.79%  "
$ws.Range("E34").Value = "  +2.72%  "
$ws.Range("D35").Value = "'7.10"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +9.63%  "
$ws.Range("D36").Value = "'0.0624"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.00%  "
$ws.Range("D37").Value = "'2.45"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.58%  "
$ws.Range("D38").Value = "'3.57"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.02%  "
$ws.Range("D39").Value = "'1.00"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("E40").Value = "  +2.79%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'103.97"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.83%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").Value = "'18.12"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.46%  "
$ws.Range("D43").Value = "1.532.55"
$ws.Range("E43").Value = "  -0.63%  "
$ws.Range("E44").Value = "  +5.96%  "
$ws.Range("E45").Value = "  +6.99%  "
$ws.Range("E46").Value = "  +0.78%  "
$ws.Range("B47").Value = "HuobiToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D47").Value = "'2.81"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.38%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "'7.81"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.74%  "
$ws.Range("B49").Value = "FTXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D49").Value = "'4.25"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.71%  "
$ws.Range("D50").Value = "2.373.48"
$ws.Range("E50").Value = "  +3.59%  "
$ws.Range("D51").Value = "'2.98"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.27%  "
